# Apply the update described by the commit diff:
# - Orders sheet: append 20 new rows (72-91) of order line items
# - Summary sheet: update G2 running tally of F-column quantities

$wb = $excel.ActiveWorkbook
$ordersWs = $wb.Worksheets.Item("Orders")

$ordersWs.Cells.Item(72, 3).Value = "624_多丁白_undefined_undefined_1bunch"
$ordersWs.Cells.Item(72, 6).Value = "'15"

$ordersWs.Cells.Item(73, 1).Value = "'10"
$ordersWs.Cells.Item(73, 3).Value = "138_卡罗拉_Carola_Rosa rugosa Thunb._20stems"
$ordersWs.Cells.Item(73, 6).Value = "'4"

$ordersWs.Cells.Item(74, 3).Value = "154_莫泊_Moab_Rosa rugosa Thunb._20stems"
$ordersWs.Cells.Item(74, 6).Value = "'11"

$ordersWs.Cells.Item(75, 3).Value = "148_坦尼克_Tineke_Rosa rugosa Thunb._20stems"
$ordersWs.Cells.Item(75, 6).Value = "'5"

$ordersWs.Cells.Item(76, 3).Value = "157_流沙_Quicksand_Rosa rugosa Thunb._20stems"
$ordersWs.Cells.Item(76, 6).Value = "'6"

$ordersWs.Cells.Item(77, 3).Value = "137_凯瑟琳_Catherine_Rosa rugosa Thunb._20stems"
$ordersWs.Cells.Item(77, 6).Value = "'11"

$ordersWs.Cells.Item(78, 3).Value = "600_康乃馨复古红_vintage red_undefined_20stems"
$ordersWs.Cells.Item(78, 6).Value = "'13"

$ordersWs.Cells.Item(79, 1).Value = "'11"
$ordersWs.Cells.Item(79, 3).Value = "135_甜蜜曼塔_sweet menta_Rosa rugosa Thunb._20stems"
$ordersWs.Cells.Item(79, 6).Value = "'8"

$ordersWs.Cells.Item(80, 3).Value = "157_流沙_Quicksand_Rosa rugosa Thunb._20stems"
$ordersWs.Cells.Item(80, 6).Value = "'7"

$ordersWs.Cells.Item(81, 3).Value = "152_白荔枝_White Ohara_Rosa rugosa Thunb._20stems"
$ordersWs.Cells.Item(81, 6).Value = "'8"

$ordersWs.Cells.Item(82, 3).Value = "203_佛罗伊德_Floyd_Rosa rugosa Thunb._20stems"
$ordersWs.Cells.Item(82, 6).Value = "'6"

$ordersWs.Cells.Item(83, 3).Value = "192_粉荔枝_Pink Ohara_Rosa rugosa Thunb._20stems"
$ordersWs.Cells.Item(83, 6).Value = "'5"

$ordersWs.Cells.Item(84, 3).Value = "147_娜欧米_Red Naomi_Rosa rugosa Thunb._20stems"
$ordersWs.Cells.Item(84, 6).Value = "'7"

$ordersWs.Cells.Item(85, 3).Value = "600_康乃馨复古红_vintage red_undefined_20stems"
$ordersWs.Cells.Item(85, 6).Value = "'7"

$ordersWs.Cells.Item(86, 3).Value = "604_康乃馨粉佳人_pink_undefined_20stems"
$ordersWs.Cells.Item(86, 6).Value = "'5"

$ordersWs.Cells.Item(87, 1).Value = "'12"
$ordersWs.Cells.Item(87, 3).Value = "192_粉荔枝_Pink Ohara_Rosa rugosa Thunb._20stems"
$ordersWs.Cells.Item(87, 6).Value = "'12"

$ordersWs.Cells.Item(88, 3).Value = "277_草莓杏仁饼_undefined_Rosa rugosa Thunb._10stems"
$ordersWs.Cells.Item(88, 6).Value = "'5"

$ordersWs.Cells.Item(89, 3).Value = "221_朱丽叶塔_Julieta_Rosa rugosa Thunb._10stems"
$ordersWs.Cells.Item(89, 6).Value = "'5"

$ordersWs.Cells.Item(90, 3).Value = "238_苏菲宝贝_undefined_Rosa rugosa Thunb._10stems"
$ordersWs.Cells.Item(90, 6).Value = "'5"

$ordersWs.Cells.Item(91, 3).Value = "244_繁星_undefined_Rosa rugosa Thunb._10stems"

$summaryWs = $wb.Worksheets.Item("Summary")
$summaryWs.Cells.Item(2, 7).Value = "'055155552510652566555525321515822555510555551255156558101576510612610551051510555510158105151051541156111387865775125550"
